$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 50
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = "2021-04-01 04:15:40.364422"
$ws.Range("D5").Value = 25
$ws.Range("E5").Value = 25
$ws.Range("F5").Value = 25
$ws.Range("G5").Value = 21
$ws.Range("H5").Value = 25
$ws.Range("I5").Value = 28

# Row 6
$ws.Range("A6").Value = 50
$ws.Range("B6").Value = 10
$ws.Range("C6").Value = "2021-04-01 04:18:00.339001"
$ws.Range("D6").Value = 19
$ws.Range("E6").Value = 19
$ws.Range("F6").Value = 19
$ws.Range("G6").Value = 18
$ws.Range("H6").Value = 19
$ws.Range("I6").Value = 18

# Row 7
$ws.Range("A7").Value = 50
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = "2021-04-01 04:20:35.999069"
$ws.Range("D7").Value = 22
$ws.Range("E7").Value = 21
$ws.Range("F7").Value = 21
$ws.Range("G7").Value = 21
$ws.Range("H7").Value = 21
$ws.Range("I7").Value = 19

# Row 8
$ws.Range("A8").Value = 50
$ws.Range("B8").Value = 10
$ws.Range("C8").Value = "2021-04-01 04:31:17.323771"
$ws.Range("D8").Value = 17
$ws.Range("E8").Value = 18
$ws.Range("F8").Value = 18
$ws.Range("G8").Value = 9
$ws.Range("H8").Value = 18
$ws.Range("I8").Value = 20

# Row 9
$ws.Range("A9").Value = 50
$ws.Range("B9").Value = 10
$ws.Range("C9").Value = "2021-04-01 04:32:29.436554"
$ws.Range("D9").Value = 27
$ws.Range("E9").Value = 25
$ws.Range("F9").Value = 26
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 26
$ws.Range("I9").Value = 26

# Row 10 (note: no value in column I for this row)
$ws.Range("A10").Value = 50
$ws.Range("B10").Value = 10
$ws.Range("C10").Value = "2021-04-04 04:51:56.682948"
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = -1
$ws.Range("F10").Value = -1
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
